# "Removed Test Case Inter-Dependency"
#
# 1. The test-case identifier stored in NewLoanInput!B2 is renamed from
#    "...VALIDATE-RANGE" to "...VALIDATE-RANGE-1st" so this scenario no
#    longer shares its name/data with another dependent test case.
# 2. The workbook's active sheet/selection is switched from
#    "Edit Repayment Schedule" back to "NewLoanInput" (cell B2), which is
#    where the renamed value now lives.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("NewLoanInput")

# Update the (previously shared / inter-dependent) test case name.
$ws1.Range("B2").Value = "2601-RBI-EI-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-VALIDATE-RANGE-1st"

# Make "NewLoanInput" the active sheet again, with B2 selected.
$ws1.Activate()
$ws1.Range("B2").Select()
